$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 34 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A34").Value = "Technische storing"
$logs.Range("B34").Value = "mailmind.test@zohomail.eu"
$logs.Range("C34").Value = "De website werkt niet goed. Is hier iets mis mee?"
$logs.Range("D34").Value = "IT / Technisch probleem"
$logs.Range("E34").Value = "Beste klant,`nBedankt voor je bericht. Om het probleem beter te kunnen onderzoeken, hebben we meer informatie nodig. Zou je alsjeblieft kunnen aangeven welke specifieke problemen je ervaart wanneer je de website probeert te gebruiken? Bijvoorbeeld, krijg je een foutmelding te zien of lukt het niet om in te loggen? Met deze details kunnen we het probleem gericht aanpakken en een oplossing bieden.`nWe horen graag meer van je, zodat we je verder kunnen helpen.`nMet vriendelijke groet,  `n[Tekstschrijver]  `nE-mailassistent"
$logs.Range("F34").Value = "2025-06-22 19:04:33"
$logs.Range("G34").Value = "Ja"

# --- Extend conditional formatting ranges to include the new row ---
$fcD = $logs.Range("D2:D33").FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($logs.Range("D2:D34"))

$fcG = $logs.Range("G2:G33").FormatConditions.Item(1)
$fcG.ModifyAppliesToRange($logs.Range("G2:G34"))

# --- "Dashboard" sheet: "IT / Technisch probleem" count 4 -> 5 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 5
